$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "_____" separator marker in column B of every user-story header row
# (rows whose column A holds the story title: 3, 8, 14, 22, 28, 32).
$ws.Range("B3").Value = "_____"
$ws.Range("B8").Value = "_____"
$ws.Range("B14").Value = "_____"
$ws.Range("B22").Value = "_____"
$ws.Range("B28").Value = "_____"
$ws.Range("B32").Value = "_____"

# Row 38's story uses the shorter "____" marker instead.
$ws.Range("B38").Value = "____"

# New user story block appended at the bottom: header row (A42) followed by
# a single task row (43) with estimate/actual time.
$ws.Range("A42").Value = "____"
$ws.Range("B43").Value = "New System parameters + admin"
$ws.Range("C43").Value = "10min"
$ws.Range("D43").Value = "15min"

# Update the view state to match where the author left the selection/scroll.
$ws.Range("D39").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
